$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 248.94444
$ws.Range("I2").Value = 93
$ws.Range("J2").Value = 794.75
$ws.Range("K2").Value = 93
$ws.Range("L2").Value = 794.75
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = -1020.75
$ws.Range("H6").Value = 2794.8
$ws.Range("I6").Value = 3244.75
$ws.Range("J6").Value = 995
$ws.Range("K6").Value = 9734.25
$ws.Range("L6").Value = 2985
$ws.Range("M6").Value = -9622.25
$ws.Range("N6").Value = -3209
$ws.Range("H21").Value = 23998.5
$ws.Range("J21").Value = 23998.5
$ws.Range("L21").Value = 23998.5
$ws.Range("N21").Value = -24934.5
$ws.Range("H23").Value = 23998.5
$ws.Range("J23").Value = 23998.5
$ws.Range("L23").Value = 23998.5
$ws.Range("N23").Value = -24466.5
$ws.Range("H28").Value = 2946.8572
$ws.Range("I28").Value = 3364.3333
$ws.Range("K28").Value = 3364.3333
$ws.Range("M28").Value = -2879.3333
$ws.Range("H29").Value = 168
$ws.Range("I29").Value = 168
$ws.Range("K29").Value = 504
$ws.Range("M29").Value = -223
$ws.Range("H80").Value = 150013500
$ws.Range("I80").Value = 500000500
$ws.Range("J80").Value = 33351168
$ws.Range("K80").Value = 1500001500
$ws.Range("L80").Value = 100053504
$ws.Range("M80").Value = -1500000502
$ws.Range("N80").Value = -100055500
$ws.Range("H83").Value = 150013500
$ws.Range("I83").Value = 500000500
$ws.Range("J83").Value = 33351168
$ws.Range("K83").Value = 4500004500
$ws.Range("L83").Value = 300160512
$ws.Range("M83").Value = -4499999508
$ws.Range("N83").Value = -300170496
$ws.Range("H87").Value = 52000
$ws.Range("J87").Value = 52000
$ws.Range("L87").Value = 52000
$ws.Range("N87").Value = -54496
$ws.Range("H90").Value = 52000
$ws.Range("J90").Value = 52000
$ws.Range("L90").Value = 156000
$ws.Range("N90").Value = -168480
$ws.Range("H98").Value = 1378.8096
$ws.Range("I98").Value = 1433.8158
$ws.Range("K98").Value = 1433.8158
$ws.Range("M98").Value = 64.18419999999992
$ws.Range("H107").Value = 2552.182
$ws.Range("I107").Value = 2726.8
$ws.Range("J107").Value = 806
$ws.Range("K107").Value = 2726.8
$ws.Range("L107").Value = 806
$ws.Range("M107").Value = -806.8000000000002
$ws.Range("N107").Value = -4646
$ws.Range("H122").Value = 1378.8096
$ws.Range("I122").Value = 1433.8158
$ws.Range("K122").Value = 4301.4474
$ws.Range("M122").Value = -1851.4474
$ws.Range("H125").Value = 1972.1666
$ws.Range("I125").Value = 2475
$ws.Range("J125").Value = 966.5
$ws.Range("K125").Value = 22275
$ws.Range("L125").Value = 8698.5
$ws.Range("M125").Value = -19815
$ws.Range("N125").Value = -13618.5
$ws.Range("H132").Value = 3206.9673
$ws.Range("I132").Value = 3362.1667
$ws.Range("J132").Value = 2633.923
$ws.Range("K132").Value = 10086.5001
$ws.Range("L132").Value = 7901.768999999999
$ws.Range("M132").Value = -7556.500100000001
$ws.Range("N132").Value = -12961.769
$ws.Range("H138").Value = 3465.1633
$ws.Range("J138").Value = 8646.362999999999
$ws.Range("L138").Value = 25939.089
$ws.Range("N138").Value = -36219.089
$ws.Range("H141").Value = 8631
$ws.Range("I141").Value = 8997
$ws.Range("K141").Value = 26991
$ws.Range("M141").Value = -21811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 4150
$ws.Range("J25").Value = 1000
$ws.Range("L25").Value = 1000
$ws.Range("N25").Value = -1804
$ws.Range("H74").Value = 1648
$ws.Range("I74").Value = 1576.963
$ws.Range("K74").Value = 1576.963
$ws.Range("M74").Value = -702.963
$ws.Range("H77").Value = 1648
$ws.Range("I77").Value = 1576.963
$ws.Range("K77").Value = 7884.815
$ws.Range("M77").Value = -3516.815
$ws.Range("H132").Value = 2540.8765
$ws.Range("I132").Value = 1744.56
$ws.Range("K132").Value = 5233.68
$ws.Range("M132").Value = -2703.68
$ws.Range("H133").Value = 116826.375
$ws.Range("J133").Value = 117642.86
$ws.Range("L133").Value = 117642.86
$ws.Range("N133").Value = -122702.86
$ws.Range("H135").Value = 85214.5
$ws.Range("J135").Value = 95429
$ws.Range("L135").Value = 95429
$ws.Range("N135").Value = -105569
$ws.Range("H137").Value = 139888.56
$ws.Range("J137").Value = 161800
$ws.Range("L137").Value = 161800
$ws.Range("N137").Value = -172000
$ws.Range("H138").Value = 87355.75
$ws.Range("J138").Value = 74712
$ws.Range("L138").Value = 74712
$ws.Range("N138").Value = -84992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15129.906
$ws.Range("I20").Value = 20702.863
$ws.Range("J20").Value = 2869.4
$ws.Range("K20").Value = 20702.863
$ws.Range("L20").Value = 2869.4
$ws.Range("M20").Value = -20455.863
$ws.Range("N20").Value = -3363.4
$ws.Range("H86").Value = 6753.1177
$ws.Range("I86").Value = 6822.778
$ws.Range("J86").Value = 6674.75
$ws.Range("K86").Value = 6822.778
$ws.Range("L86").Value = 6674.75
$ws.Range("M86").Value = -5699.778
$ws.Range("N86").Value = -8920.75
$ws.Range("H89").Value = 6753.1177
$ws.Range("I89").Value = 6822.778
$ws.Range("J89").Value = 6674.75
$ws.Range("K89").Value = 34113.89
$ws.Range("L89").Value = 33373.75
$ws.Range("M89").Value = -28497.89
$ws.Range("N89").Value = -44605.75
$ws.Range("H94").Value = 1291.2
$ws.Range("I94").Value = 1368.3334
$ws.Range("J94").Value = 1059.8
$ws.Range("K94").Value = 1368.3334
$ws.Range("L94").Value = 1059.8
$ws.Range("M94").Value = -917.3334
$ws.Range("N94").Value = -1961.8
$ws.Range("H107").Value = 3300
$ws.Range("I107").Value = 1833.3334
$ws.Range("J107").Value = 5500
$ws.Range("K107").Value = 1833.3334
$ws.Range("L107").Value = 5500
$ws.Range("M107").Value = 86.66660000000002
$ws.Range("N107").Value = -9340
$ws.Range("H115").Value = 115000.664
$ws.Range("J115").Value = 115000.664
$ws.Range("L115").Value = 115000.664
$ws.Range("N115").Value = -118134.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18717.092
$ws.Range("I60").Value = 9544.444
$ws.Range("K60").Value = 9544.444
$ws.Range("M60").Value = -9033.444
$ws.Range("H134").Value = 1880.6364
$ws.Range("I134").Value = 1401.3948
$ws.Range("K134").Value = 4204.1844
$ws.Range("M134").Value = -1669.1844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 55557556
$ws.Range("J43").Value = 55557556
$ws.Range("L43").Value = 166672668
$ws.Range("N43").Value = -166672896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4523
$ws.Range("I70").Value = 4578.2
$ws.Range("K70").Value = 4578.2
$ws.Range("M70").Value = -4308.2
$ws.Range("H73").Value = 4523
$ws.Range("I73").Value = 4578.2
$ws.Range("K73").Value = 4578.2
$ws.Range("M73").Value = -3642.2
$ws.Range("H97").Value = 737.73914
$ws.Range("I97").Value = 596.5714
$ws.Range("J97").Value = 2220
$ws.Range("K97").Value = 596.5714
$ws.Range("L97").Value = 2220
$ws.Range("M97").Value = -100.5714
$ws.Range("N97").Value = -3212
$ws.Range("H126").Value = 2132.25
$ws.Range("I126").Value = 1954.1111
$ws.Range("J126").Value = 2666.6667
$ws.Range("K126").Value = 5862.3333
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("M126").Value = -3392.3333
$ws.Range("N126").Value = -12940.0001
$ws.Range("H140").Value = 170462.33
$ws.Range("J140").Value = 164555
$ws.Range("L140").Value = 164555
$ws.Range("N140").Value = -174915

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3108.3333
$ws.Range("H40").Value = 6451.7856
$ws.Range("I40").Value = 7619.5557
$ws.Range("K40").Value = 7619.5557
$ws.Range("M40").Value = -7483.5557
$ws.Range("H46").Value = 3804.8333
$ws.Range("J46").Value = 5234.4
$ws.Range("L46").Value = 5234.4
$ws.Range("N46").Value = -5610.4
$ws.Range("H68").Value = 5997.5
$ws.Range("I68").Value = 5995
$ws.Range("J68").Value = 5998.3335
$ws.Range("K68").Value = 5995
$ws.Range("L68").Value = 5998.3335
$ws.Range("M68").Value = -5246
$ws.Range("N68").Value = -7496.3335
$ws.Range("H71").Value = 5997.5
$ws.Range("I71").Value = 5995
$ws.Range("J71").Value = 5998.3335
$ws.Range("K71").Value = 29975
$ws.Range("L71").Value = 29991.6675
$ws.Range("M71").Value = -26231
$ws.Range("N71").Value = -37479.6675
$ws.Range("H126").Value = 3108.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6772.6
$ws.Range("I126").Value = 6197.923
$ws.Range("K126").Value = 18593.769
$ws.Range("M126").Value = -16123.769
$ws.Range("H132").Value = 4417.3066
$ws.Range("I132").Value = 4037.413
$ws.Range("K132").Value = 12112.239
$ws.Range("M132").Value = -9582.239
$ws.Range("H135").Value = 97712.664
$ws.Range("J135").Value = 97712.664
$ws.Range("L135").Value = 97712.664
$ws.Range("N135").Value = -107852.664

Write-Output "Applied all cell updates"